$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Preserve the cell's existing style while forcing the new value to be
    # stored as text (avoids Excel auto-converting numeric-looking strings
    # such as "131.29" into a floating point number).
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "60.566.90"
Set-TextValue $ws.Range("E2") "  +0.50%  "
Set-TextValue $ws.Range("D3") "2.337.57"
Set-TextValue $ws.Range("E3") "  -0.17%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "551.54"
Set-TextValue $ws.Range("E5") "  +1.09%  "
Set-TextValue $ws.Range("D6") "131.29"
Set-TextValue $ws.Range("E6") "  -0.48%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.581"
Set-TextValue $ws.Range("E8") "  -0.73%  "
Set-TextValue $ws.Range("D9") "2.336.24"
Set-TextValue $ws.Range("E9") "  -0.10%  "
Set-TextValue $ws.Range("E10") "  +1.17%  "
Set-TextValue $ws.Range("D11") "5.61"
Set-TextValue $ws.Range("E11") "  +1.74%  "
Set-TextValue $ws.Range("E12") "  -0.46%  "
Set-TextValue $ws.Range("D13") "0.338"
Set-TextValue $ws.Range("E13") "  +1.26%  "
Set-TextValue $ws.Range("E14") "  +0.50%  "
Set-TextValue $ws.Range("D15") "2.754.62"
Set-TextValue $ws.Range("E15") "  -0.10%  "
Set-TextValue $ws.Range("D16") "60.491.85"
Set-TextValue $ws.Range("E16") "  +0.48%  "
Set-TextValue $ws.Range("E17") "  +1.15%  "
Set-TextValue $ws.Range("D18") "2.327.13"
Set-TextValue $ws.Range("E18") "  -0.51%  "
Set-TextValue $ws.Range("E19") "  +0.82%  "
Set-TextValue $ws.Range("E20") "  -0.92%  "
Set-TextValue $ws.Range("D21") "315.18"
Set-TextValue $ws.Range("E21") "  +0.45%  "
Set-TextValue $ws.Range("E22") "  -2.68%  "
Set-TextValue $ws.Range("E24") "  +1.28%  "
Set-TextValue $ws.Range("E25") "  -0.58%  "
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.07%  "
Set-TextValue $ws.Range("D27") "7.99"
Set-TextValue $ws.Range("E27") "  +1.06%  "
Set-TextValue $ws.Range("D28") "1.41"
Set-TextValue $ws.Range("E28") "  +3.40%  "
Set-TextValue $ws.Range("D29") "1.28"
Set-TextValue $ws.Range("E29") "  +8.99%  "
Set-TextValue $ws.Range("E30") "  -0.05%  "
Set-TextValue $ws.Range("D31") "171.19"
Set-TextValue $ws.Range("E31") "  -0.24%  "
Set-TextValue $ws.Range("D32") "0.0₃0736"
Set-TextValue $ws.Range("E32") "  +0.82%  "
Set-TextValue $ws.Range("E33") "  +2.49%  "
Set-TextValue $ws.Range("E34") "  +0.75%  "
Set-TextValue $ws.Range("E35") "  -1.34%  "
Set-TextValue $ws.Range("D36") "18.07"
Set-TextValue $ws.Range("E36") "  +0.24%  "
Set-TextValue $ws.Range("E37") "  -0.03%  "
Set-TextValue $ws.Range("E38") "  -0.05%  "
Set-TextValue $ws.Range("D39") "4.14"
Set-TextValue $ws.Range("E39") "  -0.24%  "
Set-TextValue $ws.Range("D40") "330.51"
Set-TextValue $ws.Range("E40") "  +2.70%  "
Set-TextValue $ws.Range("E41") "  +0.35%  "
Set-TextValue $ws.Range("D42") "38.06"
Set-TextValue $ws.Range("E42") "  -0.17%  "
Set-TextValue $ws.Range("D43") "138.25"
Set-TextValue $ws.Range("E43") "  -2.20%  "
Set-TextValue $ws.Range("D44") "3.52"
Set-TextValue $ws.Range("E44") "  +1.73%  "
Set-TextValue $ws.Range("E45") "  +0.69%  "
Set-TextValue $ws.Range("D46") "19.34"
Set-TextValue $ws.Range("E46") "  -0.95%  "
Set-TextValue $ws.Range("D48") "0.0499"
Set-TextValue $ws.Range("E48") "  +0.62%  "
Set-TextValue $ws.Range("D49") "0.0₆0224"
Set-TextValue $ws.Range("E49") "  +6.93%  "
Set-TextValue $ws.Range("E51") "  -0.69%  "
